$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Dhh"
$ws.Range("C2").Value2 = "Boc"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 2.264243666666667
$ws.Range("H2").Value2 = 6.792731
$ws.Range("I2").Value2 = 0.4114976873616865
$ws.Range("J2").Value2 = 0.4114976873616865
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.956431333333333
$ws.Range("N2").Value2 = 8.869294
$ws.Range("O2").Value2 = 0.06161326453142597
$ws.Range("P2").Value2 = 0.06161326453142598
$ws.Range("Q2").Value2 = 6.694080922434889
$ws.Range("R2").Value2 = 60.246728301914
$ws.Range("S2").Value2 = 0.02535371586548561
$ws.Range("T2").Value2 = 0.02535371586548562

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Dhh"
$ws.Range("C3").Value2 = "Boc"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 2.264243666666667
$ws.Range("H3").Value2 = 6.792731
$ws.Range("I3").Value2 = 0.4114976873616865
$ws.Range("J3").Value2 = 0.4114976873616865
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 40.76140833333334
$ws.Range("N3").Value2 = 122.284225
$ws.Range("O3").Value2 = 0.8494847845776018
$ws.Range("P3").Value2 = 0.8494847845776018
$ws.Range("Q3").Value2 = 92.29376066316391
$ws.Range("R3").Value2 = 830.6438459684751
$ws.Range("S3").Value2 = 0.3495610243026236
$ws.Range("T3").Value2 = 0.3495610243026236

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Dhh"
$ws.Range("C4").Value2 = "Boc"
$ws.Range("D4").Value2 = "M2"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 2.264243666666667
$ws.Range("H4").Value2 = 6.792731
$ws.Range("I4").Value2 = 0.4114976873616865
$ws.Range("J4").Value2 = 0.4114976873616865
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.007557333333333333
$ws.Range("N4").Value2 = 0.022672
$ws.Range("O4").Value2 = 0.0001574979850094596
$ws.Range("P4").Value2 = 0.0001574979850094596
$ws.Range("Q4").Value2 = 0.01711164413688889
$ws.Range("R4").Value2 = 0.154004797232
$ws.Range("S4").Value2 = 0.00006481005659551819
$ws.Range("T4").Value2 = 0.00006481005659551819

# Row 5
$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "Dhh"
$ws.Range("C5").Value2 = "Boc"
$ws.Range("D5").Value2 = "sCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 2.264243666666667
$ws.Range("H5").Value2 = 6.792731
$ws.Range("I5").Value2 = 0.4114976873616865
$ws.Range("J5").Value2 = 0.4114976873616865
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 4.258285666666667
$ws.Range("N5").Value2 = 12.774857
$ws.Range("O5").Value2 = 0.08874445290596285
$ws.Range("P5").Value2 = 0.08874445290596285
$ws.Range("Q5").Value2 = 9.641796351607445
$ws.Range("R5").Value2 = 86.776167164467
$ws.Range("S5").Value2 = 0.03651813713698181
$ws.Range("T5").Value2 = 0.03651813713698181

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Dhh"
$ws.Range("C6").Value2 = "Boc"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 2.089228666666667
$ws.Range("H6").Value2 = 6.267686
$ws.Range("I6").Value2 = 0.3796909216792509
$ws.Range("J6").Value2 = 0.3796909216792509
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.956431333333333
$ws.Range("N6").Value2 = 8.869294
$ws.Range("O6").Value2 = 0.06161326453142597
$ws.Range("P6").Value2 = 0.06161326453142598
$ws.Range("Q6").Value2 = 6.176661092631555
$ws.Range("R6").Value2 = 55.589949833684
$ws.Range("S6").Value2 = 0.02339399719760463
$ws.Range("T6").Value2 = 0.02339399719760463

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Dhh"
$ws.Range("C7").Value2 = "Boc"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 2.089228666666667
$ws.Range("H7").Value2 = 6.267686
$ws.Range("I7").Value2 = 0.3796909216792509
$ws.Range("J7").Value2 = 0.3796909216792509
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 40.76140833333334
$ws.Range("N7").Value2 = 122.284225
$ws.Range("O7").Value2 = 0.8494847845776018
$ws.Range("P7").Value2 = 0.8494847845776018
$ws.Range("Q7").Value2 = 85.15990278370558
$ws.Range("R7").Value2 = 766.4391250533502
$ws.Range("S7").Value2 = 0.3225416608087695
$ws.Range("T7").Value2 = 0.3225416608087695

# Row 8
$ws.Range("A8").Value2 = "FAPs"
$ws.Range("B8").Value2 = "Dhh"
$ws.Range("C8").Value2 = "Boc"
$ws.Range("D8").Value2 = "M2"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 2.089228666666667
$ws.Range("H8").Value2 = 6.267686
$ws.Range("I8").Value2 = 0.3796909216792509
$ws.Range("J8").Value2 = 0.3796909216792509
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.007557333333333333
$ws.Range("N8").Value2 = 0.022672
$ws.Range("O8").Value2 = 0.0001574979850094596
$ws.Range("P8").Value2 = 0.0001574979850094596
$ws.Range("Q8").Value2 = 0.01578899744355556
$ws.Range("R8").Value2 = 0.142100976992
$ws.Range("S8").Value2 = 0.00005980055509086655
$ws.Range("T8").Value2 = 0.00005980055509086655

# Row 9
$ws.Range("A9").Value2 = "FAPs"
$ws.Range("B9").Value2 = "Dhh"
$ws.Range("C9").Value2 = "Boc"
$ws.Range("D9").Value2 = "sCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 2.089228666666667
$ws.Range("H9").Value2 = 6.267686
$ws.Range("I9").Value2 = 0.3796909216792509
$ws.Range("J9").Value2 = 0.3796909216792509
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 4.258285666666667
$ws.Range("N9").Value2 = 12.774857
$ws.Range("O9").Value2 = 0.08874445290596285
$ws.Range("P9").Value2 = 0.08874445290596285
$ws.Range("Q9").Value2 = 8.896532485655777
$ws.Range("R9").Value2 = 80.068792370902
$ws.Range("S9").Value2 = 0.03369546311778591
$ws.Range("T9").Value2 = 0.03369546311778591

# Row 10
$ws.Range("A10").Value2 = "M2"
$ws.Range("B10").Value2 = "Dhh"
$ws.Range("C10").Value2 = "Boc"
$ws.Range("D10").Value2 = "ECs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.1905406666666667
$ws.Range("H10").Value2 = 0.571622
$ws.Range("I10").Value2 = 0.03462835949856721
$ws.Range("J10").Value2 = 0.03462835949856721
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 2.956431333333333
$ws.Range("N10").Value2 = 8.869294
$ws.Range("O10").Value2 = 0.06161326453142597
$ws.Range("P10").Value2 = 0.06161326453142598
$ws.Range("Q10").Value2 = 0.5633203972075556
$ws.Range("R10").Value2 = 5.069883574867999
$ws.Range("S10").Value2 = 0.002133566274074539
$ws.Range("T10").Value2 = 0.002133566274074539

# Row 11
$ws.Range("A11").Value2 = "M2"
$ws.Range("B11").Value2 = "Dhh"
$ws.Range("C11").Value2 = "Boc"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 0.1905406666666667
$ws.Range("H11").Value2 = 0.571622
$ws.Range("I11").Value2 = 0.03462835949856721
$ws.Range("J11").Value2 = 0.03462835949856721
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 40.76140833333334
$ws.Range("N11").Value2 = 122.284225
$ws.Range("O11").Value2 = 0.8494847845776018
$ws.Range("P11").Value2 = 0.8494847845776018
$ws.Range("Q11").Value2 = 7.766705918105557
$ws.Range("R11").Value2 = 69.90035326295
$ws.Range("S11").Value2 = 0.02941626450891612
$ws.Range("T11").Value2 = 0.02941626450891612

# Row 12
$ws.Range("A12").Value2 = "M2"
$ws.Range("B12").Value2 = "Dhh"
$ws.Range("C12").Value2 = "Boc"
$ws.Range("D12").Value2 = "M2"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 0.1905406666666667
$ws.Range("H12").Value2 = 0.571622
$ws.Range("I12").Value2 = 0.03462835949856721
$ws.Range("J12").Value2 = 0.03462835949856721
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.007557333333333333
$ws.Range("N12").Value2 = 0.022672
$ws.Range("O12").Value2 = 0.0001574979850094596
$ws.Range("P12").Value2 = 0.0001574979850094596
$ws.Range("Q12").Value2 = 0.001439979331555555
$ws.Range("R12").Value2 = 0.012959813984
$ws.Range("S12").Value2 = 0.000005453896845207516
$ws.Range("T12").Value2 = 0.000005453896845207516

# Row 13
$ws.Range("A13").Value2 = "M2"
$ws.Range("B13").Value2 = "Dhh"
$ws.Range("C13").Value2 = "Boc"
$ws.Range("D13").Value2 = "sCs"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 0.1905406666666667
$ws.Range("H13").Value2 = 0.571622
$ws.Range("I13").Value2 = 0.03462835949856721
$ws.Range("J13").Value2 = 0.03462835949856721
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 4.258285666666667
$ws.Range("N13").Value2 = 12.774857
$ws.Range("O13").Value2 = 0.08874445290596285
$ws.Range("P13").Value2 = 0.08874445290596285
$ws.Range("Q13").Value2 = 0.8113765897837778
$ws.Range("R13").Value2 = 7.302389308053999
$ws.Range("S13").Value2 = 0.003073074818731349
$ws.Range("T13").Value2 = 0.003073074818731349

# Row 14
$ws.Range("A14").Value2 = "sCs"
$ws.Range("B14").Value2 = "Dhh"
$ws.Range("C14").Value2 = "Boc"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 0.9584326666666668
$ws.Range("H14").Value2 = 2.875298
$ws.Range("I14").Value2 = 0.1741830314604954
$ws.Range("J14").Value2 = 0.1741830314604954
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 2.956431333333333
$ws.Range("N14").Value2 = 8.869294
$ws.Range("O14").Value2 = 0.06161326453142597
$ws.Range("P14").Value2 = 0.06161326453142598
$ws.Range("Q14").Value2 = 2.833540366623556
$ws.Range("R14").Value2 = 25.501863299612
$ws.Range("S14").Value2 = 0.0107319851942612
$ws.Range("T14").Value2 = 0.0107319851942612

# Row 15
$ws.Range("A15").Value2 = "sCs"
$ws.Range("B15").Value2 = "Dhh"
$ws.Range("C15").Value2 = "Boc"
$ws.Range("D15").Value2 = "FAPs"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 0.9584326666666668
$ws.Range("H15").Value2 = 2.875298
$ws.Range("I15").Value2 = 0.1741830314604954
$ws.Range("J15").Value2 = 0.1741830314604954
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 40.76140833333334
$ws.Range("N15").Value2 = 122.284225
$ws.Range("O15").Value2 = 0.8494847845776018
$ws.Range("P15").Value2 = 0.8494847845776018
$ws.Range("Q15").Value2 = 39.06706528600557
$ws.Range("R15").Value2 = 351.6035875740501
$ws.Range("S15").Value2 = 0.1479658349572926
$ws.Range("T15").Value2 = 0.1479658349572926

# Row 16
$ws.Range("A16").Value2 = "sCs"
$ws.Range("B16").Value2 = "Dhh"
$ws.Range("C16").Value2 = "Boc"
$ws.Range("D16").Value2 = "M2"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 0.9584326666666668
$ws.Range("H16").Value2 = 2.875298
$ws.Range("I16").Value2 = 0.1741830314604954
$ws.Range("J16").Value2 = 0.1741830314604954
$ws.Range("K16").Value2 = 1
$ws.Range("L16").Value2 = 0.3333333333333333
$ws.Range("M16").Value2 = 0.007557333333333333
$ws.Range("N16").Value2 = 0.022672
$ws.Range("O16").Value2 = 0.0001574979850094596
$ws.Range("P16").Value2 = 0.0001574979850094596
$ws.Range("Q16").Value2 = 0.007243195139555557
$ws.Range("R16").Value2 = 0.06518875625600001
$ws.Range("S16").Value2 = 0.00002743347647786734
$ws.Range("T16").Value2 = 0.00002743347647786734

# Row 17
$ws.Range("A17").Value2 = "sCs"
$ws.Range("B17").Value2 = "Dhh"
$ws.Range("C17").Value2 = "Boc"
$ws.Range("D17").Value2 = "sCs"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 0.9584326666666668
$ws.Range("H17").Value2 = 2.875298
$ws.Range("I17").Value2 = 0.1741830314604954
$ws.Range("J17").Value2 = 0.1741830314604954
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 4.258285666666667
$ws.Range("N17").Value2 = 12.774857
$ws.Range("O17").Value2 = 0.08874445290596285
$ws.Range("P17").Value2 = 0.08874445290596285
$ws.Range("Q17").Value2 = 4.081280086931779
$ws.Range("R17").Value2 = 36.731520782386
$ws.Range("S17").Value2 = 0.01545777783246378
$ws.Range("T17").Value2 = 0.01545777783246378
